$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-11: columns E-M get new benchmark values ---
$updatedEM = @{
    2 = @(-9.84, 72.54, 66.08, 14.59, 62.71, 73.42, -10.78, 83.59, 75.44)
    3 = @(13.07, 57.28, 66, 13.57, 60.38, 70, 24.93, 58.94, 78.4)
    4 = @(10.7, 54.19, 60.68, 20.97, 53.27, 67.44, 20.7, 56.38, 71.04)
    5 = @(-3.51, 60.67, 58.63, 8.17, 63.41, 69.1, -3.31, 71.82, 69.53)
    6 = @(9.65, 43.94, 48.6, 20.67, 47.37, 59.7, 39.16, 36.94, 60.7)
    7 = @(21.5, 44.18, 56.27, 18.48, 54.27, 66.56, 67.77, 21.54, 66.82)
    8 = @(30.13, 41.33, 59.26, 11.99, 59.64, 67.59, 43.85, 41.79, 74.38)
    9 = @(3.46, 57.01, 58.96, 17.56, 55.51, 67.35, 28.84, 48.99, 68.84)
    10 = @(56.55, 20.36, 46.91, 75.31, 14.69, 59.54, 76.03, 14.59, 60.87)
    11 = @(-16.74, 62.08, 53.16, 22.82, 50.65, 65.6, 24.36, 49.71, 65.71)
}
foreach ($r in $updatedEM.Keys) {
    $vals = $updatedEM[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $vals[$i]
    }
}

# --- Append new rows 12-21 (columns A-M) for the new experiment runs ---
$newRows = @{
    12 = @(4.06, "Planar", 2.5, 0, 1.31, 378.98, 383.74, 0.03, 378.35, 378.33, 6.3, 379.13, 404.43)
    13 = @(2.5, "Planar", 0, 2.5, 2.65, 604.82, 621.2, 0.1, 610.28, 610.8, 62.2, 247.01, 653.6)
    14 = @(4.55, "Planar", 2.5, 2.5, 64.85, 134.41, 381.98, 0.15, 377.9, 378.02, 64.62, 141.95, 400.88)
    15 = @(15.37, "Planar", 10, 0, 72.06, 61.57, 220.36, 0.08, 222.22, 222.45, 65.75, 77.68, 226.81)
    16 = @(10, "Planar", 0, 10, 28.95, 133.22, 187.5, 0.88, 186.84, 188.5, 32.96, 132.33, 197.4)
    17 = @(19.26, "Planar", 10, 10, -0.57, 196.52, 195.38, 0.18, 197.88, 198.23, 55.49, 89.43, 200.88)
    18 = @(45.63, "Gradual", 0, 2.5, 43.8, 84.73, 150.76, 2.99, 152.42, 157.11, 45.99, 85.03, 157.44)
    19 = @(43.22, "Gradual", 2.5, 2.5, 35.4, 106.43, 164.74, 1.47, 167.13, 169.6, 37.64, 106.63, 170.99)
    20 = @(81, "Gradual", 0, 10, 37.94, 77.79, 125.33, 42.09, 77.56, 133.94, 4.69, 127.79, 134.07)
    21 = @(85.86, "Gradual", 10, 10, 57.95, 53.95, 128.3, 0.99, 135.66, 137.03, 3.12, 133.4, 137.7)
}
foreach ($r in $newRows.Keys) {
    $vals = $newRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 1 + $i).Value = $vals[$i]
    }
}
